$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 232-233, shifting the existing rows 232:257 down to 234:259.
$ws.Rows("232:233").Insert()

# Row 232: new weekly record, "1a amarillo"
$ws.Cells.Item(232, 1).Value = 4
$ws.Cells.Item(232, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(232, 3).Value = "Los Lagos"
$ws.Cells.Item(232, 4).Value = 44449
$ws.Cells.Item(232, 5).Value = 10
$ws.Cells.Item(232, 6).Value = "Fruta"
$ws.Cells.Item(232, 7).Value = 100102
$ws.Cells.Item(232, 8).Value = "Cítricos"
$ws.Cells.Item(232, 9).Value = 100102003
$ws.Cells.Item(232, 10).Value = "Limón"
$ws.Cells.Item(232, 11).Value = "Sin especificar"
$ws.Cells.Item(232, 12).Value = "1a amarillo"
$ws.Cells.Item(232, 13).Value = 1000
$ws.Cells.Item(232, 14).Value = 8000
$ws.Cells.Item(232, 15).Value = 8000
$ws.Cells.Item(232, 16).Value = 8000
$ws.Cells.Item(232, 17).Value = "`$/malla 16 kilos"
$ws.Cells.Item(232, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(232, 19).Value = 500
$ws.Cells.Item(232, 20).Value = 16

# Row 233: new weekly record, "2a amarillo"
$ws.Cells.Item(233, 1).Value = 4
$ws.Cells.Item(233, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(233, 3).Value = "Los Lagos"
$ws.Cells.Item(233, 4).Value = 44449
$ws.Cells.Item(233, 5).Value = 10
$ws.Cells.Item(233, 6).Value = "Fruta"
$ws.Cells.Item(233, 7).Value = 100102
$ws.Cells.Item(233, 8).Value = "Cítricos"
$ws.Cells.Item(233, 9).Value = 100102003
$ws.Cells.Item(233, 10).Value = "Limón"
$ws.Cells.Item(233, 11).Value = "Sin especificar"
$ws.Cells.Item(233, 12).Value = "2a amarillo"
$ws.Cells.Item(233, 13).Value = 500
$ws.Cells.Item(233, 14).Value = 7000
$ws.Cells.Item(233, 15).Value = 7000
$ws.Cells.Item(233, 16).Value = 7000
$ws.Cells.Item(233, 17).Value = "`$/malla 16 kilos"
$ws.Cells.Item(233, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(233, 19).Value = 438
$ws.Cells.Item(233, 20).Value = 16

# Make sure the date cells keep the original date/time number format.
$ws.Range("D232:D233").NumberFormat = "YYYY-MM-DD HH:MM:SS"
